$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Beijing_Enviornment" worksheet between
#    "Beijing_Nearest_Neighbor" (1st) and "Beijing_1ring_Neighbors" (was 2nd).
# ---------------------------------------------------------------------------
$wsNearest = $wb.Worksheets.Item(1)
$wsEnv = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNearest)
$wsEnv.Name = "Beijing_Enviornment"

# ---------------------------------------------------------------------------
# 2. Populate the new sheet with the environmental land-class table.
# ---------------------------------------------------------------------------
$rows = @(
    @("station_id", "terrain", "landclass"),
    @("beijing_grid_202", "flat", "suburbs"),
    @("beijing_grid_203", "flat", "suburbs"),
    @("beijing_grid_203", "flat", "farm"),
    @("beijing_grid_204", "mountain", "forest"),
    @("beijing_grid_223", "hills", "forest"),
    @("beijing_grid_224", "mountain", "forest"),
    @("beijing_grid_224", "flat", "farm"),
    @("beijing_grid_225", "flat", "park"),
    @("beijing_grid_282", "flat", "park"),
    @("beijing_grid_283", "flat", "city"),
    @("beijing_grid_303", "flat", "city"),
    @("beijing_grid_304", "flat", "park"),
    @("chaoyang_meo", "flat", "park"),
    @("hadian_meo", "flat", "city"),
    @("aotizhongxin_aq", "flat", "park"),
    @("beibuxinqu_aq", "flat", "park"),
    @("yanqin_aq", "mountains", "farm"),
    @("badaling_aq", "flat", "park"),
    @("yanqin_meo", "flat", "park"),
    @("yanqin_meo", $null, $null)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $val = $rowVals[$c]
        if ($val -ne $null) {
            $wsEnv.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# Column widths to match the sibling sheets (A:D).
$wsEnv.Range("A1:D1").ColumnWidth = 20.26953125

# ---------------------------------------------------------------------------
# 3. Register the (stale/hidden) AutoFilter defined name for the new sheet
#    without leaving a visible <autoFilter> element, matching the source
#    workbook's state.
# ---------------------------------------------------------------------------
$wsEnv.Names.Add("_xlnm._FilterDatabase", "=Beijing_Enviornment!`$A`$1:`$H`$659") | Out-Null
$envFilterName = $wb.Names.Item($wb.Names.Count)
$envFilterName.Visible = $false

# ---------------------------------------------------------------------------
# 4. Beijing_Nearest_Neighbor (1st sheet): selection now A1:B36, no longer
#    the tab shown when the workbook opens.
# ---------------------------------------------------------------------------
$wsNearest.Range("A1:B36").Select()

# ---------------------------------------------------------------------------
# 5. Beijing_1ring_Neighbors (last sheet): select entire grid and turn the
#    AutoFilter on for its data range.
# ---------------------------------------------------------------------------
$ws1Ring = $wb.Worksheets.Item(3)
$ws1Ring.Cells.Select()
$ws1Ring.Range("A1:F725").AutoFilter() | Out-Null

# ---------------------------------------------------------------------------
# 6. Make "Beijing_Enviornment" the active tab/sheet (matches activeTab="1").
# ---------------------------------------------------------------------------
$wsEnv.Range("A1:C21").Select()
$wsEnv.Activate()
